# LoginTestDataProvider.xlsx edit:
#  - sheet1 "loginTest1": row3 col B value changes from "admin123" to "admin";
#    row4 (the "field30"/"field40" row) is removed.
#  - sheet2 "loginTest2": row3 (the "login test2 field 10" row) is removed.
#  - Selections / active sheet tweaked to match the post-edit screenshot state.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("loginTest1")
$ws2 = $wb.Worksheets.Item("loginTest2")

# --- sheet1: loginTest1 ---
# Row 3, column B: "admin123" -> "admin"
$ws1.Cells.Item(3, 2).Value = "admin"

# Remove row 4 entirely (login test1 field30 / field40)
$ws1.Rows.Item(4).Delete() | Out-Null

# --- sheet2: loginTest2 ---
# Remove row 3 entirely (login test2 field 10)
$ws2.Rows.Item(3).Delete() | Out-Null

# --- selections / active sheet ---
$ws1.Select()
$ws1.Range("C10").Select()

$ws2.Range("C10").Select()

$ws1.Activate()
